{"js": "// \"study instance\" -> \"SOP instance\" (x2) and \"series/study\" -> \"series/SOP\"\n// (x1): the author renamed \"study instance unique identifier\" to \"SOP\n// instance unique identifier\" everywhere it appears in the Abstract.\n//\n// The source document has track changes turned on; temporarily turn it\n// off so the edit lands as plain text (not a tracked insertion/\n// deletion), then restore the document's original track-changes\n// setting so that's not an unintended side effect of this edit.\ncontext.document.changeTrackingMode = Word.ChangeTrackingMode.off;\nawait context.sync();\n\nconst body = context.document.body;\n\n// The \"Validation Methods\" paragraph originally has its sentence spread\n// across many same-formatted runs (Word split them apart while the\n// authors were editing, tracked via now-meaningless w:rsidR bookkeeping).\n// Re-set each stretch of text (the parts between the existing grammar-\n// check markers) as a single piece so the runs collapse back together,\n// matching how Word normalizes a paragraph it has just edited.\nconst beforeGramStart = body.search(\n  \"Validation Methods: Publicly available brain MRI and TCI lung 4\",\n  { matchCase: true }\n);\nbeforeGramStart.load(\"items\");\nawait context.sync();\nif (beforeGramStart.items.length > 0) {\n  beforeGramStart.items[0].insertText(\n    \"Validation Methods: Publicly available brain MRI and TCI lung 4\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nconst betweenGramMarks = body.search(\"DCT  images\", { matchCase: true });\nbetweenGramMarks.load(\"items\");\nawait context.sync();\nif (betweenGramMarks.items.length > 0) {\n  betweenGramMarks.items[0].insertText(\"DCT  images\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// This stretch also contains one of the \"study instance\" -> \"SOP\n// instance\" renames.\nconst oldTail =\n  \" were used to evaluate the software. The ability to change the frame of reference, series instance identifier, and study instance identifier using the program was evaluated with both the \";\nconst newTail =\n  \" were used to evaluate the software. The ability to change the frame of reference, series instance identifier, and SOP instance identifier using the program was evaluated with both the \";\nconst tail = body.search(oldTail, { matchCase: true });\ntail.load(\"items\");\nawait context.sync();\nif (tail.items.length > 0) {\n  tail.items[0].insertText(newTail, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// The other two \"study instance\" -> \"SOP instance\" / \"series/study\" ->\n// \"series/SOP\" renames sit in their own, already-single-run paragraphs.\nconst studyInstance = body.search(\"study instance\", { matchCase: true });\nstudyInstance.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < studyInstance.items.length; i++) {\n  studyInstance.items[i].insertText(\"SOP instance\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst seriesStudy = body.search(\"series/study\", { matchCase: true });\nseriesStudy.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < seriesStudy.items.length; i++) {\n  seriesStudy.items[i].insertText(\"series/SOP\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Restore track changes to the state it was in originally.\ncontext.document.changeTrackingMode = Word.ChangeTrackingMode.trackAll;\nawait context.sync();\n", "ps1": "# \"study instance\" -> \"SOP instance\" (x2) and \"series/study\" -> \"series/SOP\"\n# (x1): the author renamed \"study instance unique identifier\" to \"SOP\n# instance unique identifier\" everywhere it appears in the Abstract.\n#\n# The source document has track changes turned on; temporarily turn it\n# off so the edit lands as plain text (not a tracked insertion/\n# deletion), then restore the document's original track-changes setting\n# so that's not an unintended side effect of this edit.\n$d = $word.ActiveDocument\n$d.TrackRevisions = $false\n\n# The \"Validation Methods\" paragraph originally has its sentence spread\n# across many same-formatted runs. Replace each stretch of text (the\n# parts between the existing grammar-check markers) as a single piece so\n# the runs collapse back together, matching how Word normalizes a\n# paragraph it has just edited.\n$beforeGramStart = \"Validation Methods: Publicly available brain MRI and TCI lung 4\"\n$rng1 = $d.Content\n$rng1.Find.Execute($beforeGramStart, $false, $false, $false, $false, $false, $true, 1, $false, $beforeGramStart, 2)\n\n$betweenGramMarks = \"DCT  images\"\n$rng2 = $d.Content\n$rng2.Find.Execute($betweenGramMarks, $false, $false, $false, $false, $false, $true, 1, $false, $betweenGramMarks, 2)\n\n# This stretch also contains one of the \"study instance\" -> \"SOP\n# instance\" renames.\n$oldTail = \" were used to evaluate the software. The ability to change the frame of reference, series instance identifier, and study instance identifier using the program was evaluated with both the \"\n$newTail = \" were used to evaluate the software. The ability to change the frame of reference, series instance identifier, and SOP instance identifier using the program was evaluated with both the \"\n$rng3 = $d.Content\n$rng3.Find.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)\n\n# The other two \"study instance\" -> \"SOP instance\" / \"series/study\" ->\n# \"series/SOP\" renames sit in their own, already-single-run paragraphs.\n$rng4 = $d.Content\n$rng4.Find.Execute(\"study instance\", $false, $false, $false, $false, $false, $true, 1, $false, \"SOP instance\", 2)\n\n$rng5 = $d.Content\n$rng5.Find.Execute(\"series/study\", $false, $false, $false, $false, $false, $true, 1, $false, \"series/SOP\", 2)\n\n# Restore track changes to the state it was in originally.\n$d.TrackRevisions = $true\n"}
